# EPBDS-12991: add(CustomDatatype[][], CustomDatatype[], CustomDatatype[])
# returns Object[][] instead of CustomDatatype[][].
#
# The fix changed the expected test results for Step2/Step3 of the
# "Test mySpr" table on Sheet1: the lengths computed by addAll(...) now
# come back as 7 instead of 5 (E17/F17), and the empty placeholder cell
# F14 is no longer populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# F14 was a completely empty placeholder cell - drop it entirely.
$ws.Range("F14").ClearContents()

# Updated expected values for _res_.$Step2 / _res_.$Step3 results.
$ws.Range("E17").Value = 7
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 1
